$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1961538461538462
$ws.Range("C2").Value = 0.573076923076923
$ws.Range("J2").Value = 0.02692307692307692
$ws.Range("P2").Value = 0.1307692307692308
$ws.Range("S2").Value = 0.07307692307692308
$ws.Range("B3").Value = 0.0130718954248366
$ws.Range("C3").Value = 0.0392156862745098
$ws.Range("J3").Value = 0.0130718954248366
$ws.Range("P3").Value = 0.7516339869281046
$ws.Range("S3").Value = 0.1830065359477124
$ws.Range("J4").Value = 0.08108108108108109
$ws.Range("P4").Value = 0.6486486486486487
$ws.Range("S4").Value = 0.2702702702702703
$ws.Range("B6").Value = 0.06435643564356436
$ws.Range("D6").Value = 0.01485148514851485
$ws.Range("E6").Value = 0.004950495049504951
$ws.Range("F6").Value = 0.0396039603960396
$ws.Range("J6").Value = 0.2227722772277228
$ws.Range("O6").Value = 0.02475247524752475
$ws.Range("Q6").Value = 0.1336633663366337
$ws.Range("R6").Value = 0.103960396039604
$ws.Range("S6").Value = 0.3910891089108911
$ws.Range("B7").Value = 0.1153846153846154
$ws.Range("D7").Value = 0.03846153846153846
$ws.Range("F7").Value = 0.04395604395604396
$ws.Range("J7").Value = 0.1483516483516484
$ws.Range("O7").Value = 0.03296703296703297
$ws.Range("Q7").Value = 0.1318681318681319
$ws.Range("R7").Value = 0.06593406593406594
$ws.Range("S7").Value = 0.4230769230769231
$ws.Range("B8").Value = 0.07936507936507936
$ws.Range("D8").Value = 0.01388888888888889
$ws.Range("F8").Value = 0.08134920634920635
$ws.Range("J8").Value = 0.125
$ws.Range("O8").Value = 0.02182539682539682
$ws.Range("Q8").Value = 0.1726190476190476
$ws.Range("R8").Value = 0.07539682539682539
$ws.Range("S8").Value = 0.4305555555555556
$ws.Range("B9").Value = 0.07729468599033816
$ws.Range("D9").Value = 0.01449275362318841
$ws.Range("F9").Value = 0.06763285024154589
$ws.Range("J9").Value = 0.106280193236715
$ws.Range("O9").Value = 0.01932367149758454
$ws.Range("Q9").Value = 0.1352657004830918
$ws.Range("R9").Value = 0.106280193236715
$ws.Range("S9").Value = 0.4734299516908212
$ws.Range("B10").Value = 0.09076923076923077
$ws.Range("D10").Value = 0.01307692307692308
$ws.Range("E10").Value = 0.0007692307692307692
$ws.Range("F10").Value = 0.06
$ws.Range("J10").Value = 0.1046153846153846
$ws.Range("O10").Value = 0.01153846153846154
$ws.Range("Q10").Value = 0.1961538461538462
$ws.Range("R10").Value = 0.09384615384615384
$ws.Range("S10").Value = 0.4292307692307692
$ws.Range("G11").Value = 0.1402214022140221
$ws.Range("J11").Value = 0.1070110701107011
$ws.Range("K11").Value = 0.1808118081180812
$ws.Range("L11").Value = 0.5645756457564576
$ws.Range("S11").Value = 0.007380073800738007
$ws.Range("G12").Value = 0.7658227848101266
$ws.Range("J12").Value = 0.1835443037974684
$ws.Range("L12").Value = 0.04430379746835443
$ws.Range("S12").Value = 0.006329113924050633
$ws.Range("G13").Value = 0.6829268292682927
$ws.Range("J13").Value = 0.3170731707317073
$ws.Range("F15").Value = 0.0091324200913242
$ws.Range("H15").Value = 0.2328767123287671
$ws.Range("I15").Value = 0.0821917808219178
$ws.Range("J15").Value = 0.3515981735159817
$ws.Range("K15").Value = 0.0502283105022831
$ws.Range("M15").Value = 0.0182648401826484
$ws.Range("N15").Value = 0.0045662100456621
$ws.Range("O15").Value = 0.0365296803652968
$ws.Range("S15").Value = 0.2146118721461187
$ws.Range("F16").Value = 0.01176470588235294
$ws.Range("H16").Value = 0.2294117647058823
$ws.Range("I16").Value = 0.08235294117647059
$ws.Range("J16").Value = 0.4294117647058823
$ws.Range("K16").Value = 0.09411764705882353
$ws.Range("M16").Value = 0.01176470588235294
$ws.Range("N16").Value = 0.005882352941176471
$ws.Range("O16").Value = 0.05882352941176471
$ws.Range("S16").Value = 0.07647058823529412
$ws.Range("F17").Value = 0.00477326968973747
$ws.Range("H17").Value = 0.2195704057279236
$ws.Range("I17").Value = 0.1073985680190931
$ws.Range("J17").Value = 0.4367541766109785
$ws.Range("K17").Value = 0.09785202863961814
$ws.Range("M17").Value = 0.007159904534606206
$ws.Range("N17").Value = 0.002386634844868735
$ws.Range("O17").Value = 0.06205250596658711
$ws.Range("S17").Value = 0.06205250596658711
$ws.Range("F18").Value = 0.009302325581395349
$ws.Range("H18").Value = 0.2325581395348837
$ws.Range("I18").Value = 0.1069767441860465
$ws.Range("J18").Value = 0.4186046511627907
$ws.Range("K18").Value = 0.05581395348837209
$ws.Range("M18").Value = 0.05116279069767442
$ws.Range("O18").Value = 0.06976744186046512
$ws.Range("S18").Value = 0.05581395348837209
$ws.Range("F19").Value = 0.01256873527101336
$ws.Range("H19").Value = 0.2144540455616654
$ws.Range("I19").Value = 0.08483896307934014
$ws.Range("J19").Value = 0.4092694422623723
$ws.Range("K19").Value = 0.1068342498036135
$ws.Range("M19").Value = 0.01649646504320503
$ws.Range("N19").Value = 0.001571091908876669
$ws.Range("O19").Value = 0.07305577376276512
$ws.Range("S19").Value = 0.08091123330714847
